# Auto-applied price/volume updates for cryptos.xlsx
# Generated from the commit diff (Fri Feb  9 10:48:04 UTC 2024 GitHub Actions run)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "46.566.62"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.96%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.464.63"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.65%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.29"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.07"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.10%  "
$ws.Range("E7").Value = "  +1.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  +0.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.06"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.62%  "
$ws.Range("E11").Value = "  +1.64%  "
$ws.Range("E12").Value = "  +0.68%  "
$ws.Range("E13").Value = "  -2.89%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.07"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.94%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.848.87"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.64%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.457.17"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.84%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.842"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.07%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "46.460.71"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.62"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.46%  "
$ws.Range("E20").Value = "  +1.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0935"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.83%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "70.42"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "248.74"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.65%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.38"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.41%  "
$ws.Range("E25").Value = "  +2.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.08"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.25%  "
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("E28").Value = "  +0.98%  "
$ws.Range("E29").Value = "  +3.23%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.04"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.84%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "49.52"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.130"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.56"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.22%  "
$ws.Range("E34").Value = "  +3.13%  "
$ws.Range("E35").Value = "  +0.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0767"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.98%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.62"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.38%  "
$ws.Range("E38").Value = "  +0.86%  "
$ws.Range("E39").Value = "  +2.49%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "123.31"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.45%  "
$ws.Range("E41").Value = "  +1.83%  "
$ws.Range("E42").Value = "  +0.93%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.69"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.64%  "
$ws.Range("E44").Value = "  +0.73%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.982.18"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.09%  "
$ws.Range("E46").Value = "  +0.70%  "
$ws.Range("E47").Value = "  -2.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.80"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.83%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.93"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.84%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.34"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +15.63%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.22"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.10%  "
